$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 134, shifting existing rows 134:189 down to 135:190.
$ws.Rows(134).Insert()

# Populate the newly inserted row 134 with its data (same master/category fields as
# the row that used to occupy 134, but with its own date and price figures).
$ws.Cells.Item(134, 1).Value = 11
$ws.Cells.Item(134, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(134, 3).Value = "Bíobío"
$ws.Cells.Item(134, 4).Value = 44992
$ws.Cells.Item(134, 5).Value = 8
$ws.Cells.Item(134, 6).Value = 100112043
$ws.Cells.Item(134, 7).Value = "Pepino ensalada"
$ws.Cells.Item(134, 8).Value = "Sin especificar"
$ws.Cells.Item(134, 9).Value = "Primera"
$ws.Cells.Item(134, 10).Value = 100
$ws.Cells.Item(134, 11).Value = 8000
$ws.Cells.Item(134, 12).Value = 8500
$ws.Cells.Item(134, 13).Value = 8250
$ws.Cells.Item(134, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(134, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(134, 16).Value = 138
$ws.Cells.Item(134, 17).Value = 60
$ws.Cells.Item(134, 18).Value = "Hortaliza"
